# Insert a new leading column (index "#") before the existing data.
# This shifts the existing "name"/"birthday"/"random_int" columns one
# position to the right (A->B, B->C, C->D) and keeps their custom
# widths, number formats and the bold header style intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Insert()

# Header for the new column, styled the same as the other bold headers.
$ws.Range("A1").Value = "#"
$ws.Range("A1").Font.Bold = $true

# Row index values for the new column.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
